$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 61.13720810128422
$ws.Range("C2").Value = 95.31067699233704
$ws.Range("D2").Value = 99.33273066436571
$ws.Range("E2").Value = 98.94492273052317
$ws.Range("F2").Value = 98.36697756451264
$ws.Range("G2").Value = 97.4076759688394
$ws.Range("H2").Value = 96.03454094752266
$ws.Range("B3").Value = 70.26108722209689
$ws.Range("C3").Value = 95.3163587664579
$ws.Range("D3").Value = 99.80924244665123
$ws.Range("E3").Value = 99.04053924570317
$ws.Range("F3").Value = 98.57715559428648
$ws.Range("G3").Value = 97.57300415138546
$ws.Range("H3").Value = 96.08113580816774
$ws.Range("B4").Value = 82.62678759805662
$ws.Range("C4").Value = 94.9414740935602
$ws.Range("D4").Value = 99.20938114632808
$ws.Range("E4").Value = 98.71189271957832
$ws.Range("F4").Value = 98.41846724632862
$ws.Range("G4").Value = 97.44913118947842
$ws.Range("H4").Value = 96.0110550149187
$ws.Range("B5").Value = 76.02576500087805
$ws.Range("C5").Value = 94.90069951064287
$ws.Range("D5").Value = 99.27287226637506
$ws.Range("E5").Value = 98.88827420406025
$ws.Range("F5").Value = 98.35810753054093
$ws.Range("G5").Value = 97.38975000588403
$ws.Range("H5").Value = 95.96731985759092
$ws.Range("B6").Value = 72.55875255079552
$ws.Range("C6").Value = 95.38332739227181
$ws.Range("D6").Value = 99.34559258537894
$ws.Range("E6").Value = 98.8684690848535
$ws.Range("F6").Value = 98.33511628117778
$ws.Range("G6").Value = 97.47987514698821
$ws.Range("H6").Value = 96.01776895712371
